# Update the loginTest data sheet so the 3rd/4th data rows use the same
# credentials ("satya") as row 2 instead of "james"/"test", and leave the
# sheet's selection on the data block (A2:B4) rather than the old
# whole-column selection at C1:C1048576.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginTest")

$ws.Range("A3").Value = "satya"
$ws.Range("B3").Value = "satya"
$ws.Range("A4").Value = "satya"
$ws.Range("B4").Value = "satya"

$ws.Range("A2:B4").Select()
